$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.202382356775502
$ws.Range("C2").Value = 0.3362148038159205
$ws.Range("D2").Value = 0.02038419427572791
$ws.Range("F2").Value = 0.3834850062737729
$ws.Range("G2").Value = 0.2346391638160554
$ws.Range("H2").Value = 0.403514836126341
$ws.Range("I2").Value = 0.2472595070819583

$ws.Range("B3").Value = 1.053040124645179
$ws.Range("C3").Value = 0.2941235210392961
$ws.Range("D3").Value = 0.01801444738944724
$ws.Range("F3").Value = 0.3852676611739483
$ws.Range("G3").Value = 0.2371158725544404
$ws.Range("H3").Value = 0.4107019864815342
$ws.Range("I3").Value = 0.2566918009935053

$ws.Range("B4").Value = 0.9610690578164736
$ws.Range("C4").Value = 0.268170013557949
$ws.Range("D4").Value = 0.0165544840562788
$ws.Range("F4").Value = 0.3869363474800025
$ws.Range("G4").Value = 0.2391348160629647
$ws.Range("H4").Value = 0.4155451368749326
$ws.Range("I4").Value = 0.2628852308347849

$ws.Range("B5").Value = 0.9235230683864302
$ws.Range("C5").Value = 0.257566726401933
$ws.Range("D5").Value = 0.01595833916409362
$ws.Range("F5").Value = 0.3877599089007333
$ws.Range("G5").Value = 0.2400818410016754
$ws.Range("H5").Value = 0.4176265868993241
$ws.Range("I5").Value = 0.2655096991807744

$ws.Range("B6").Value = 0.9172845948604049
$ws.Range("C6").Value = 0.2558044402767905
$ws.Range("D6").Value = 0.0158592786410594
$ws.Range("F6").Value = 0.3879053074084311
$ws.Range("G6").Value = 0.240246571301995
$ws.Range("H6").Value = 0.4179787122153797
$ws.Range("I6").Value = 0.2659515517492217

$ws.Range("B7").Value = 0.9605629682777703
$ws.Range("C7").Value = 0.2680271227123399
$ws.Range("D7").Value = 0.01654644903251068
$ws.Range("F7").Value = 0.3869468741331232
$ws.Range("G7").Value = 0.2391470860523128
$ws.Range("H7").Value = 0.4155727719213544
$ws.Range("I7").Value = 0.2629202186655171

$ws.Range("B8").Value = 1.150947208842183
$ws.Range("C8").Value = 0.3217246736274433
$ws.Range("D8").Value = 0.01956815253590349
$ws.Range("F8").Value = 0.3839800167516714
$ws.Range("G8").Value = 0.2353891257462521
$ws.Range("H8").Value = 0.4059034634468901
$ws.Range("I8").Value = 0.2504280464055748

$ws.Range("B9").Value = 1.52204682681122
$ws.Range("C9").Value = 0.4261433249340598
$ws.Range("D9").Value = 0.02545308584232231
$ws.Range("F9").Value = 0.3827544952042246
$ws.Range("G9").Value = 0.2320179726888512
$ws.Range("H9").Value = 0.3903713021648301
$ws.Range("I9").Value = 0.2291414819066695

$ws.Range("B10").Value = 1.793264766615721
$ws.Range("C10").Value = 0.5023096981594222
$ws.Range("D10").Value = 0.02975042343562961
$ws.Range("F10").Value = 0.3847056611268513
$ws.Range("G10").Value = 0.2320407806703457
$ws.Range("H10").Value = 0.3810720933916798
$ws.Range("I10").Value = 0.2154879272343528

$ws.Range("B11").Value = 1.916327736768039
$ws.Range("C11").Value = 0.5368383205453142
$ws.Range("D11").Value = 0.03169938827828389
$ws.Range("F11").Value = 0.3862236418952421
$ws.Range("G11").Value = 0.2326072530207881
$ws.Range("H11").Value = 0.3773049530921782
$ws.Range("I11").Value = 0.2097138200365265

$ws.Range("B12").Value = 1.962881654335945
$ws.Range("C12").Value = 0.5498958510283956
$ws.Range("D12").Value = 0.03243652739367064
$ws.Range("F12").Value = 0.3868899965390327
$ws.Range("G12").Value = 0.2329028010644691
$ws.Range("H12").Value = 0.375945415820766
$ws.Range("I12").Value = 0.2075906653220354

$ws.Range("B13").Value = 1.95285757301076
$ws.Range("C13").Value = 0.5470844733388844
$ws.Range("D13").Value = 0.03227781159104381
$ws.Range("F13").Value = 0.3867424004138584
$ws.Range("G13").Value = 0.2328355283146237
$ws.Range("H13").Value = 0.3762352303967305
$ws.Range("I13").Value = 0.2080450974915546

$ws.Range("B14").Value = 1.920158716058154
$ws.Range("C14").Value = 0.5379129290578817
$ws.Range("D14").Value = 0.03176005125675374
$ws.Range("F14").Value = 0.3862766238987092
$ws.Range("G14").Value = 0.2326299376998548
$ws.Range("H14").Value = 0.3771917581238426
$ws.Range("I14").Value = 0.2095378736123337

$ws.Range("B15").Value = 1.900123495746811
$ws.Range("C15").Value = 0.5322927730771312
$ws.Range("D15").Value = 0.03144279079336343
$ws.Range("F15").Value = 0.3860032680884657
$ws.Range("G15").Value = 0.2325145927140397
$ws.Range("H15").Value = 0.3777863958740113
$ws.Range("I15").Value = 0.2104605113198215

$ws.Range("B16").Value = 1.785215793138434
$ws.Range("C16").Value = 0.5000507190491703
$ws.Range("D16").Value = 0.02962293144443606
$ws.Range("F16").Value = 0.384619217718793
$ws.Range("G16").Value = 0.232015042076128
$ws.Range("H16").Value = 0.3813276382750388
$ws.Range("I16").Value = 0.2158741176581014

$ws.Range("B17").Value = 1.714641416909444
$ws.Range("C17").Value = 0.4802402272178483
$ws.Range("D17").Value = 0.02850496216854026
$ws.Range("F17").Value = 0.3839322310042519
$ws.Range("G17").Value = 0.2318517765120731
$ws.Range("H17").Value = 0.3836189766942226
$ws.Range("I17").Value = 0.2193074904485073

$ws.Range("B18").Value = 1.674019329972225
$ws.Range("C18").Value = 0.4688345040759714
$ws.Range("D18").Value = 0.02786138139976657
$ws.Range("F18").Value = 0.3835963768379358
$ws.Range("G18").Value = 0.2318101490416069
$ws.Range("H18").Value = 0.3849804606661422
$ws.Range("I18").Value = 0.2213233724508665

$ws.Range("B19").Value = 1.660260371697575
$ws.Range("C19").Value = 0.4649708042574616
$ws.Range("D19").Value = 0.02764338204953987
$ws.Range("F19").Value = 0.3834928173692802
$ws.Range("G19").Value = 0.2318050009180865
$ws.Range("H19").Value = 0.3854489061399988
$ws.Range("I19").Value = 0.2220129574357674

$ws.Range("B20").Value = 1.722157254957267
$ws.Range("C20").Value = 0.482350255933568
$ws.Range("D20").Value = 0.02862402958571408
$ws.Range("F20").Value = 0.3839992200259204
$ws.Range("G20").Value = 0.2318637380087694
$ws.Range("H20").Value = 0.383370547905713
$ws.Range("I20").Value = 0.2189377447813516

$ws.Range("B21").Value = 1.929764464050152
$ws.Range("C21").Value = 0.5406073183202125
$ws.Range("D21").Value = 0.03191215453644247
$ws.Range("F21").Value = 0.386410942700472
$ws.Range("G21").Value = 0.2326881167916071
$ws.Range("H21").Value = 0.3769089812659132
$ws.Range("I21").Value = 0.2090976848168218

$ws.Range("B22").Value = 2.06517083002899
$ws.Range("C22").Value = 0.5785782368241144
$ws.Range("D22").Value = 0.03405591912168404
$ws.Range("F22").Value = 0.3885210465322047
$ws.Range("G22").Value = 0.2336997389049031
$ws.Range("H22").Value = 0.3730766699321748
$ws.Range("I22").Value = 0.2030362816591875

$ws.Range("B23").Value = 1.992927897267407
$ws.Range("C23").Value = 0.5583220629773678
$ws.Range("D23").Value = 0.03291224210093446
$ws.Range("F23").Value = 0.387345702061026
$ws.Range("G23").Value = 0.23311619854708
$ws.Range("H23").Value = 0.3750861653109325
$ws.Range("I23").Value = 0.2062373657367971

$ws.Range("B24").Value = 1.718759493753737
$ws.Range("C24").Value = 0.4813963630293756
$ws.Range("D24").Value = 0.02857020183664361
$ws.Range("F24").Value = 0.3839687502513343
$ws.Range("G24").Value = 0.2318581675858695
$ws.Range("H24").Value = 0.3834827250228727
$ws.Range("I24").Value = 0.2191047760389484

$ws.Range("B25").Value = 1.421901177709856
$ws.Range("C25").Value = 0.3979910456552034
$ws.Range("D25").Value = 0.02386556859976707
$ws.Range("F25").Value = 0.3825890069455156
$ws.Range("G25").Value = 0.2324954889619377
$ws.Range("H25").Value = 0.3942037907656157
$ws.Range("I25").Value = 0.2345534106892959
